# Reorder the roster rows 10-19 on Sheet1 so that each player again lines up
# with their correct Position / Team, but in the new row order introduced by
# the commit (a pure re-shuffle of the same 10 players/rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 1).Value = "Daniel Gafford"
$ws.Cells.Item(10, 2).Value = "PF,C"
$ws.Cells.Item(10, 3).Value = "Dallas Mavericks"

$ws.Cells.Item(11, 1).Value = "Jalen Duren"
$ws.Cells.Item(11, 2).Value = "C"
$ws.Cells.Item(11, 3).Value = "Detroit Pistons"

$ws.Cells.Item(12, 1).Value = "Keegan Murray"
$ws.Cells.Item(12, 2).Value = "SF,PF"
$ws.Cells.Item(12, 3).Value = "Sacramento Kings"

$ws.Cells.Item(13, 1).Value = "Jarrett Allen"
$ws.Cells.Item(13, 2).Value = "C"
$ws.Cells.Item(13, 3).Value = "Cleveland Cavaliers"

$ws.Cells.Item(14, 1).Value = "Karl-Anthony Towns"
$ws.Cells.Item(14, 2).Value = "PF,C"
$ws.Cells.Item(14, 3).Value = "New York Knicks"

$ws.Cells.Item(15, 1).Value = "Taurean Prince"
$ws.Cells.Item(15, 2).Value = "SG,SF"
$ws.Cells.Item(15, 3).Value = "Milwaukee Bucks"

$ws.Cells.Item(16, 1).Value = "Austin Reaves"
$ws.Cells.Item(16, 2).Value = "PG,SG"
$ws.Cells.Item(16, 3).Value = "Los Angeles Lakers"

$ws.Cells.Item(17, 1).Value = "Kevin Durant"
$ws.Cells.Item(17, 2).Value = "SF,PF"
$ws.Cells.Item(17, 3).Value = "Phoenix Suns"

$ws.Cells.Item(18, 1).Value = "Franz Wagner"
$ws.Cells.Item(18, 2).Value = "SF,PF"
$ws.Cells.Item(18, 3).Value = "Orlando Magic"

$ws.Cells.Item(19, 1).Value = "Mark Williams"
$ws.Cells.Item(19, 2).Value = "C"
$ws.Cells.Item(19, 3).Value = "Charlotte Hornets"
